$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add English translation of course name (new B4/C4 cells)
$ws.Range("B4").Value = "Multivariate Statistics"
$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("B4").Font.Bold = $false
$ws.Range("C4").Value = "Multivariate Statistics"

# Update activation date (keep as literal text, not an auto-converted date)
# Use a scratch cell formatted as Text so Excel doesn't reinterpret the
# dd/mm/yyyy-looking string as a date serial, then paste just the value
# into the real destination cells so their existing formatting is kept.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "01/01/2021"
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Update "Programa resumido" (short syllabus content)
$resumido = "Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e variância condicional/marginal. Estatística: Regressão Logística simples, Teste Qui-Quadrado, Testes de normalidade, Testes não-paramétricos. Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência"
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# Update "Programa" (full syllabus content)
$programa = "Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e Variância condicional/marginal.Estatística: Regressão Logística simples (coeficiente de associação, sensitividade e especificidade, risco relativo, razão de chances), Teste Qui-Quadrado (testes de aderência, homogeneidade e independência), Teste de normalidade (Shapiro-Wilk, Teste de Kolmogorov-Smirnov), Testes não-paramétricos para amostras pareadas e independentes.Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Update "Método" (grading method)
$metodo = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: N = (N1+...+Nn)/n"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Update "Critério" (grading criterion)
$criterio = "NF≥ 5,0."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Update "Norma de recuperação" (recovery norm)
$recuperacao = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# Update "Bibliografia" (bibliography)
$bibliografia = "G.C. Runger, D. Montgomery. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009. D. C. Montgomery, E. A. Peck, G. G. Vining, Introduction to Linear Regression Analysis, 4th ed., Hoboken: John Wiley, 2006.W. J. Conover, Practical Nonparametric Statistics, 3rd ed., New York: John Wiley d Sons, 1999.R. A. Johnson, D. W. Wichern, Applied Multivariate Statistical Analysis, 6th ed., New Jersey: Prentice Hall, 2007."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
